$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1238.2609
$ws.Range("I38").Value = 89.36364
$ws.Range("J38").Value = 2291.4167
$ws.Range("K38").Value = 268.09092
$ws.Range("L38").Value = 6874.250100000001
$ws.Range("M38").Value = 103.90908
$ws.Range("N38").Value = -7618.250100000001
$ws.Range("H58").Value = 813.6316
$ws.Range("J58").Value = 1265.875
$ws.Range("L58").Value = 3797.625
$ws.Range("N58").Value = -4097.625
$ws.Range("H86").Value = 2875.7273
$ws.Range("I86").Value = 3019
$ws.Range("J86").Value = 2625
$ws.Range("K86").Value = 3019
$ws.Range("L86").Value = 2625
$ws.Range("M86").Value = -1896
$ws.Range("N86").Value = -4871
$ws.Range("H89").Value = 2875.7273
$ws.Range("I89").Value = 3019
$ws.Range("J89").Value = 2625
$ws.Range("K89").Value = 15095
$ws.Range("L89").Value = 13125
$ws.Range("M89").Value = -9479
$ws.Range("N89").Value = -24357
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H131").Value = 730
$ws.Range("I131").Value = 730
$ws.Range("K131").Value = 2190
$ws.Range("M131").Value = 2850

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2721.4
$ws.Range("I88").Value = 1950
$ws.Range("J88").Value = 3235.6667
$ws.Range("K88").Value = 1950
$ws.Range("L88").Value = 3235.6667
$ws.Range("M88").Value = -1544
$ws.Range("N88").Value = -4047.6667
$ws.Range("H91").Value = 2721.4
$ws.Range("I91").Value = 1950
$ws.Range("J91").Value = 3235.6667
$ws.Range("K91").Value = 1950
$ws.Range("L91").Value = 3235.6667
$ws.Range("M91").Value = -546
$ws.Range("N91").Value = -6043.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3136.4
$ws.Range("I86").Value = 3237.75
$ws.Range("J86").Value = 2956.2222
$ws.Range("K86").Value = 3237.75
$ws.Range("L86").Value = 2956.2222
$ws.Range("M86").Value = -2114.75
$ws.Range("N86").Value = -5202.2222
$ws.Range("H89").Value = 3136.4
$ws.Range("I89").Value = 3237.75
$ws.Range("J89").Value = 2956.2222
$ws.Range("K89").Value = 16188.75
$ws.Range("L89").Value = 14781.111
$ws.Range("M89").Value = -10572.75
$ws.Range("N89").Value = -26013.111
$ws.Range("H105").Value = 76924470
$ws.Range("I105").Value = 90910190
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 90910190
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -90908443
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 845.7353000000001
$ws.Range("I107").Value = 701.9545000000001
$ws.Range("K107").Value = 701.9545000000001
$ws.Range("M107").Value = 1218.0455

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1311.5667
$ws.Range("I31").Value = 1288.5532
$ws.Range("J31").Value = 1394.7693
$ws.Range("K31").Value = 1288.5532
$ws.Range("L31").Value = 1394.7693
$ws.Range("M31").Value = -993.5532000000001
$ws.Range("N31").Value = -1984.7693
$ws.Range("H32").Value = 1310
$ws.Range("I32").Value = 1310
$ws.Range("K32").Value = 1310
$ws.Range("M32").Value = -994
$ws.Range("H34").Value = 1311.5667
$ws.Range("I34").Value = 1288.5532
$ws.Range("J34").Value = 1394.7693
$ws.Range("K34").Value = 1288.5532
$ws.Range("L34").Value = 1394.7693
$ws.Range("M34").Value = -1086.5532
$ws.Range("N34").Value = -1798.7693
$ws.Range("H99").Value = 1679.25
$ws.Range("I99").Value = 1583.4445
$ws.Range("J99").Value = 1966.6666
$ws.Range("K99").Value = 1583.4445
$ws.Range("L99").Value = 1966.6666
$ws.Range("M99").Value = -85.44450000000006
$ws.Range("N99").Value = -4962.6666
$ws.Range("H126").Value = 1679.25
$ws.Range("I126").Value = 1583.4445
$ws.Range("J126").Value = 1966.6666
$ws.Range("K126").Value = 4750.333500000001
$ws.Range("L126").Value = 5899.9998
$ws.Range("M126").Value = -2280.333500000001
$ws.Range("N126").Value = -10839.9998
$ws.Range("H132").Value = 2964.2222
$ws.Range("I132").Value = 1900
$ws.Range("J132").Value = 3097.25
$ws.Range("K132").Value = 5700
$ws.Range("L132").Value = 9291.75
$ws.Range("M132").Value = -3170
$ws.Range("N132").Value = -14351.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 11230.625
$ws.Range("I70").Value = 22318
$ws.Range("J70").Value = 6190.909
$ws.Range("K70").Value = 66954
$ws.Range("L70").Value = 18572.727
$ws.Range("M70").Value = -66639
$ws.Range("N70").Value = -19202.727
$ws.Range("H73").Value = 11230.625
$ws.Range("I73").Value = 22318
$ws.Range("J73").Value = 6190.909
$ws.Range("K73").Value = 66954
$ws.Range("L73").Value = 18572.727
$ws.Range("M73").Value = -65862
$ws.Range("N73").Value = -20756.727
$ws.Range("H74").Value = 4750
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 18000
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -20122
$ws.Range("H77").Value = 4750
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 54000
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -64608
$ws.Range("H87").Value = 1296
$ws.Range("I87").Value = 1125
$ws.Range("J87").Value = 1980
$ws.Range("K87").Value = 3375
$ws.Range("L87").Value = 5940
$ws.Range("M87").Value = -2127
$ws.Range("N87").Value = -8436
$ws.Range("H88").Value = 5644.4443
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 6250
$ws.Range("K88").Value = 2400
$ws.Range("L88").Value = 18750
$ws.Range("M88").Value = -1972
$ws.Range("N88").Value = -19606
$ws.Range("H90").Value = 1296
$ws.Range("I90").Value = 1125
$ws.Range("J90").Value = 1980
$ws.Range("K90").Value = 10125
$ws.Range("L90").Value = 17820
$ws.Range("M90").Value = -3885
$ws.Range("N90").Value = -30300
$ws.Range("H91").Value = 5644.4443
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 6250
$ws.Range("K91").Value = 2400
$ws.Range("L91").Value = 18750
$ws.Range("M91").Value = -918
$ws.Range("N91").Value = -21714

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5330
$ws.Range("I80").Value = 4816.6665
$ws.Range("J80").Value = 6100
$ws.Range("K80").Value = 4816.6665
$ws.Range("L80").Value = 6100
$ws.Range("M80").Value = -3818.6665
$ws.Range("N80").Value = -8096
$ws.Range("H83").Value = 5330
$ws.Range("I83").Value = 4816.6665
$ws.Range("J83").Value = 6100
$ws.Range("K83").Value = 24083.3325
$ws.Range("L83").Value = 30500
$ws.Range("M83").Value = -19091.3325
$ws.Range("N83").Value = -40484
$ws.Range("H102").Value = 1385.125
$ws.Range("I102").Value = 1470.1538
$ws.Range("J102").Value = 1016.6667
$ws.Range("K102").Value = 1470.1538
$ws.Range("L102").Value = 1016.6667
$ws.Range("M102").Value = 151.8462
$ws.Range("N102").Value = -4260.6667
$ws.Range("H132").Value = 3144.9312
$ws.Range("I132").Value = 3325.0625
$ws.Range("J132").Value = 2923.2307
$ws.Range("K132").Value = 9975.1875
$ws.Range("L132").Value = 8769.6921
$ws.Range("M132").Value = -7445.1875
$ws.Range("N132").Value = -13829.6921

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 2000
$ws.Range("I19").Value = 2000
$ws.Range("K19").Value = 2000
$ws.Range("M19").Value = -1830
$ws.Range("H93").Value = 1025.5
$ws.Range("J93").Value = 1550
$ws.Range("L93").Value = 1550
$ws.Range("N93").Value = -4046
$ws.Range("H132").Value = 4359.6
$ws.Range("J132").Value = 3839.6
$ws.Range("L132").Value = 11518.8
$ws.Range("N132").Value = -16578.8
$ws.Range("H135").Value = 56193
$ws.Range("J135").Value = 56193
$ws.Range("L135").Value = 56193
$ws.Range("N135").Value = -66333
$ws.Range("H136").Value = 2020.8
$ws.Range("I136").Value = 1912
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5736
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3186
$ws.Range("N136").Value = -14100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 421.55554
$ws.Range("I100").Value = 423
$ws.Range("K100").Value = 846
$ws.Range("M100").Value = -305
$ws.Range("H107").Value = 459.42856
$ws.Range("I107").Value = 352.3
$ws.Range("J107").Value = 556.8182
$ws.Range("K107").Value = 1056.9
$ws.Range("L107").Value = 1670.4546
$ws.Range("M107").Value = 863.0999999999999
$ws.Range("N107").Value = -5510.4546
$ws.Range("H132").Value = 1782.8125
$ws.Range("I132").Value = 1737.1052
$ws.Range("J132").Value = 1956.5
$ws.Range("K132").Value = 5211.3156
$ws.Range("L132").Value = 5869.5
$ws.Range("M132").Value = -2681.3156
$ws.Range("N132").Value = -10929.5
